$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B and H must hold the text "2", "2" and "05/10/2024"
# (not the numeric 2 or a date serial). Excel auto-detects numbers/dates
# from a bare string, so force text via NumberFormat, then strip the
# formatting back off so the cells keep the workbook's default style.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"

$ws.Range("A5").Value = "2"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = 2231231231232
$ws.Range("D5").Value = "A"
$ws.Range("E5").Value = "2@gmail.com"
$ws.Range("F5").Value = 22222222
$ws.Range("G5").Value = "Masculino"
$ws.Range("H5").Value = "05/10/2024"
$ws.Range("I5").Value = "NEUTRO"

$ws.Range("A5").ClearFormats()
$ws.Range("B5").ClearFormats()
$ws.Range("H5").ClearFormats()
